$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.605.90'
$ws.Range('E2').Value = '  -2.68%  '
$ws.Range('D3').Value = '3.176.39'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '596.68'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').Value = '151.90'
$ws.Range('E6').Value = '  -3.76%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '3.173.38'
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  -3.99%  '
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -4.88%  '
$ws.Range('D11').Value = '5.56'
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('D12').Value = '0.477'
$ws.Range('E12').Value = '  -5.81%  '
$ws.Range('D13').Value = '0.0000258'
$ws.Range('E13').Value = '  -5.57%  '
$ws.Range('D14').Value = '36.90'
$ws.Range('E14').Value = '  -5.30%  '
$ws.Range('D15').Value = '3.710.26'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '64.810.67'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').Value = '3.192.44'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '7.03'
$ws.Range('E19').Value = '  -5.26%  '
$ws.Range('D20').Value = '482.99'
$ws.Range('E20').Value = '  -5.25%  '
$ws.Range('D21').Value = '14.77'
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('D22').Value = '0.716'
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('D23').Value = '7.76'
$ws.Range('E23').Value = '  -3.83%  '
$ws.Range('D24').Value = '13.89'
$ws.Range('E24').Value = '  -5.75%  '
$ws.Range('D25').Value = '85.00'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Value = '2.92'
$ws.Range('E27').Value = '  -2.24%  '
$ws.Range('D28').Value = '8.67'
$ws.Range('E28').Value = '  -5.26%  '
$ws.Range('E29').Value = '  +31.04%  '
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  -5.36%  '
$ws.Range('D31').Value = '6.95'
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').Value = '2.73'
$ws.Range('E33').Value = '  -9.06%  '
$ws.Range('D34').Value = '26.88'
$ws.Range('E34').Value = '  -4.75%  '
$ws.Range('E35').Value = '  -6.40%  '
$ws.Range('D36').Value = '6.12'
$ws.Range('E36').Value = '  -6.00%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = '54.47'
$ws.Range('E37').Value = '  -3.02%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '3.27'
$ws.Range('E38').Value = '  +7.55%  '
$ws.Range('D39').Value = '472.76'
$ws.Range('E39').Value = '  -8.01%  '
$ws.Range('D40').Value = '0.0₃0729'
$ws.Range('E40').Value = '  -5.06%  '
$ws.Range('D41').Value = '0.0403'
$ws.Range('E41').Value = '  -4.18%  '
$ws.Range('D42').Value = '0.124'
$ws.Range('E42').Value = '  -3.41%  '
$ws.Range('D43').Value = '8.51'
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('D44').Value = '2.903.74'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('D45').Value = '2.43'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('E46').Value = '  -7.62%  '
$ws.Range('D47').Value = '27.35'
$ws.Range('E47').Value = '  -4.12%  '
$ws.Range('D49').Value = '2.35'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('D51').Value = '120.86'
$ws.Range('E51').Value = '  -2.05%  '
